$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns before writing so values like
# "40.545.23", "0.998", "3.81" etc. are stored as text (matching the workbook's
# original inlineStr cells) instead of being auto-coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '40.545.23'
$ws.Range("E2").Value = '  -2.62%  '
$ws.Range("D3").Value = '2.368.51'
$ws.Range("E3").Value = '  -3.99%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '310.25'
$ws.Range("E5").Value = '  -2.49%  '
$ws.Range("D6").Value = '86.59'
$ws.Range("E6").Value = '  -6.30%  '
$ws.Range("D7").Value = '0.527'
$ws.Range("E7").Value = '  -4.38%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -4.36%  '
$ws.Range("D10").Value = '0.0838'
$ws.Range("E10").Value = '  -3.65%  '
$ws.Range("D11").Value = '30.49'
$ws.Range("E11").Value = '  -7.49%  '
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '2.730.75'
$ws.Range("E13").Value = '  -4.08%  '
$ws.Range("D14").Value = '6.53'
$ws.Range("E14").Value = '  -5.12%  '
$ws.Range("D15").Value = '15.02'
$ws.Range("E15").Value = '  -3.06%  '
$ws.Range("D16").Value = '2.364.86'
$ws.Range("E16").Value = '  -4.88%  '
$ws.Range("D17").Value = '0.756'
$ws.Range("E17").Value = '  -4.89%  '
$ws.Range("D18").Value = '40.434.63'
$ws.Range("E18").Value = '  -2.76%  '
$ws.Range("D19").Value = '0.0₃0908'
$ws.Range("E19").Value = '  -4.04%  '
$ws.Range("D20").Value = '6.11'
$ws.Range("E20").Value = '  -5.19%  '
$ws.Range("D21").Value = '68.33'
$ws.Range("E21").Value = '  -3.55%  '
$ws.Range("D22").Value = '10.71'
$ws.Range("E22").Value = '  -5.03%  '
$ws.Range("D23").Value = '234.47'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("E24").Value = '  -6.32%  '
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  -7.75%  '
$ws.Range("D27").Value = '23.53'
$ws.Range("E27").Value = '  -5.36%  '
$ws.Range("D28").Value = '2.17'
$ws.Range("E28").Value = '  -3.63%  '
$ws.Range("D29").Value = '9.26'
$ws.Range("E29").Value = '  -4.64%  '
$ws.Range("D30").Value = '33.65'
$ws.Range("E30").Value = '  -8.38%  '
$ws.Range("D31").Value = '152.04'
$ws.Range("E31").Value = '  -3.68%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '5.19'
$ws.Range("E33").Value = '  -5.25%  '
$ws.Range("D34").Value = '0.0726'
$ws.Range("E34").Value = '  -4.49%  '
$ws.Range("E35").Value = '  -5.29%  '
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("D37").Value = '0.0993'
$ws.Range("E37").Value = '  -4.37%  '
$ws.Range("D38").Value = '15.76'
$ws.Range("E38").Value = '  -8.74%  '
$ws.Range("D39").Value = '2.74'
$ws.Range("E39").Value = '  -5.59%  '
$ws.Range("D40").Value = '1.70'
$ws.Range("E40").Value = '  -7.85%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '3.81'
$ws.Range("E41").Value = '  -4.47%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value = '2.39'
$ws.Range("E42").Value = '  -5.84%  '
$ws.Range("D43").Value = '1.955.04'
$ws.Range("E43").Value = '  -1.67%  '
$ws.Range("E44").Value = '  -5.32%  '
$ws.Range("D45").Value = '17.57'
$ws.Range("E45").Value = '  -7.24%  '
$ws.Range("D46").Value = '9.41'
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '2.69'
$ws.Range("E47").Value = '  -9.01%  '
$ws.Range("D48").Value = '2.589.04'
$ws.Range("E48").Value = '  -4.44%  '
$ws.Range("D49").Value = '93.00'
$ws.Range("E49").Value = '  -4.67%  '
$ws.Range("D50").Value = '72.20'
$ws.Range("E50").Value = '  -4.76%  '
$ws.Range("D51").Value = '50.16'
$ws.Range("E51").Value = '  -4.44%  '

# Restore the original (unstyled) cell style now that the values are committed as text.
$ws.Range("D2:E51").Style = "Normal"
